$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 1007.125
$ws.Range("I33").Value = 1079.5714
$ws.Range("K33").Value = 1079.5714
$ws.Range("M33").Value = -850.5714
$ws.Range("H101").Value = 1441.8334
$ws.Range("I101").Value = 1216.5
$ws.Range("J101").Value = 1892.5
$ws.Range("K101").Value = 3649.5
$ws.Range("L101").Value = 5677.5
$ws.Range("M101").Value = -2027.5
$ws.Range("N101").Value = -8921.5
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()
$ws.Range("H125").Value = 7814470.5
$ws.Range("I125").Value = 31251276
$ws.Range("J125").Value = 2202.5
$ws.Range("K125").Value = 281261484
$ws.Range("L125").Value = 19822.5
$ws.Range("M125").Value = -281259024
$ws.Range("N125").Value = -24742.5
$ws.Range("H132").Value = 2118.6667
$ws.Range("I132").Value = 2169.9348
$ws.Range("K132").Value = 6509.8044
$ws.Range("M132").Value = -3979.8044
$ws.Range("H137").Value = 3088093.5
$ws.Range("I137").Value = 10417905
$ws.Range("J137").Value = 1857.0526
$ws.Range("K137").Value = 31253715
$ws.Range("L137").Value = 5571.1578
$ws.Range("M137").Value = -31251165
$ws.Range("N137").Value = -10671.1578
$ws.Range("H138").Value = 2887.5417
$ws.Range("I138").Value = 1180.3478
$ws.Range("J138").Value = 3425.4246
$ws.Range("K138").Value = 3541.0434
$ws.Range("L138").Value = 10276.2738
$ws.Range("M138").Value = 1598.9566
$ws.Range("N138").Value = -20556.2738
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1598.5714
$ws.Range("I2").Value = 1565
$ws.Range("J2").Value = 1800
$ws.Range("K2").Value = 1565
$ws.Range("L2").Value = 1800
$ws.Range("M2").Value = -1452
$ws.Range("N2").Value = -2026
$ws.Range("H32").Value = 9626674
$ws.Range("I32").Value = 10214520
$ws.Range("K32").Value = 10214520
$ws.Range("M32").Value = -10214233
$ws.Range("H45").Value = 2169.3809
$ws.Range("I45").Value = 1746.8
$ws.Range("J45").Value = 3225.8333
$ws.Range("K45").Value = 1746.8
$ws.Range("L45").Value = 3225.8333
$ws.Range("M45").Value = -1369.8
$ws.Range("N45").Value = -3979.8333
$ws.Range("H61").Value = 15154667
$ws.Range("I61").Value = 23811692
$ws.Range("J61").Value = 4875
$ws.Range("K61").Value = 23811692
$ws.Range("L61").Value = 4875
$ws.Range("M61").Value = -23811480
$ws.Range("N61").Value = -5299
$ws.Range("H74").Value = 12823125
$ws.Range("I74").Value = 2022.5238
$ws.Range("J74").Value = 27781078
$ws.Range("K74").Value = 2022.5238
$ws.Range("L74").Value = 27781078
$ws.Range("M74").Value = -1148.5238
$ws.Range("N74").Value = -27782826
$ws.Range("H76").Value = 69950
$ws.Range("J76").Value = 69950
$ws.Range("L76").Value = 69950
$ws.Range("N76").Value = -70626
$ws.Range("H77").Value = 12823125
$ws.Range("I77").Value = 2022.5238
$ws.Range("J77").Value = 27781078
$ws.Range("K77").Value = 10112.619
$ws.Range("L77").Value = 138905390
$ws.Range("M77").Value = -5744.618999999999
$ws.Range("N77").Value = -138914126
$ws.Range("H79").Value = 69950
$ws.Range("J79").Value = 69950
$ws.Range("L79").Value = 69950
$ws.Range("N79").Value = -72290
$ws.Range("H97").Value = 1148.0952
$ws.Range("I97").Value = 1225
$ws.Range("J97").Value = 1100.7693
$ws.Range("K97").Value = 1225
$ws.Range("L97").Value = 1100.7693
$ws.Range("M97").Value = -729
$ws.Range("N97").Value = -2092.7693
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()
$ws.Range("H116").Value = 1598.5714
$ws.Range("I116").Value = 1565
$ws.Range("J116").Value = 1800
$ws.Range("K116").Value = 1565
$ws.Range("L116").Value = 1800
$ws.Range("M116").Value = 729
$ws.Range("N116").Value = -6388
$ws.Range("H136").Value = 15154667
$ws.Range("I136").Value = 23811692
$ws.Range("J136").Value = 4875
$ws.Range("K136").Value = 71435076
$ws.Range("L136").Value = 14625
$ws.Range("M136").Value = -71432526
$ws.Range("N136").Value = -19725
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1598.5714
$ws.Range("I3").Value = 1565
$ws.Range("J3").Value = 1800
$ws.Range("K3").Value = 1565
$ws.Range("L3").Value = 1800
$ws.Range("M3").Value = -1451
$ws.Range("N3").Value = -2028
$ws.Range("H22").Value = 605.8
$ws.Range("I22").Value = 933.3333
$ws.Range("J22").Value = 114.5
$ws.Range("K22").Value = 933.3333
$ws.Range("L22").Value = 114.5
$ws.Range("M22").Value = -760.3333
$ws.Range("N22").Value = -460.5
$ws.Range("H98").Value = 80000
$ws.Range("J98").Value = 80000
$ws.Range("L98").Value = 80000
$ws.Range("N98").Value = -85990
$ws.Range("H134").Value = 2854.838
$ws.Range("I134").Value = 2897.96
$ws.Range("J134").Value = 2765
$ws.Range("K134").Value = 8693.880000000001
$ws.Range("L134").Value = 8295
$ws.Range("M134").Value = -6158.880000000001
$ws.Range("N134").Value = -13365
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5186.7954
$ws.Range("I31").Value = 1706.7391
$ws.Range("J31").Value = 6520.817
$ws.Range("K31").Value = 1706.7391
$ws.Range("L31").Value = 6520.817
$ws.Range("M31").Value = -1411.7391
$ws.Range("N31").Value = -7110.817
$ws.Range("H34").Value = 5186.7954
$ws.Range("I34").Value = 1706.7391
$ws.Range("J34").Value = 6520.817
$ws.Range("K34").Value = 1706.7391
$ws.Range("L34").Value = 6520.817
$ws.Range("M34").Value = -1504.7391
$ws.Range("N34").Value = -6924.817
$ws.Range("H58").Value = 1160.0625
$ws.Range("I58").Value = 897
$ws.Range("J58").Value = 2300
$ws.Range("K58").Value = 897
$ws.Range("L58").Value = 2300
$ws.Range("M58").Value = -694
$ws.Range("N58").Value = -2706
$ws.Range("H62").Value = 3834.36
$ws.Range("I62").Value = 4429.7
$ws.Range("J62").Value = 3437.4666
$ws.Range("K62").Value = 4429.7
$ws.Range("L62").Value = 3437.4666
$ws.Range("M62").Value = -3805.7
$ws.Range("N62").Value = -4685.4666
$ws.Range("H65").Value = 3834.36
$ws.Range("I65").Value = 4429.7
$ws.Range("J65").Value = 3437.4666
$ws.Range("K65").Value = 22148.5
$ws.Range("L65").Value = 17187.333
$ws.Range("M65").Value = -19028.5
$ws.Range("N65").Value = -23427.333
$ws.Range("H81").Value = 98000
$ws.Range("J81").Value = 98000
$ws.Range("L81").Value = 98000
$ws.Range("N81").Value = -99996
$ws.Range("H84").Value = 98000
$ws.Range("J84").Value = 98000
$ws.Range("L84").Value = 294000
$ws.Range("N84").Value = -303984
$ws.Range("H134").Value = 5559524.5
$ws.Range("I134").Value = 6414178
$ws.Range("J134").Value = 4275.8335
$ws.Range("K134").Value = 19242534
$ws.Range("L134").Value = 12827.5005
$ws.Range("M134").Value = -19239999
$ws.Range("N134").Value = -17897.5005
$ws.Range("H136").Value = 1160.0625
$ws.Range("I136").Value = 897
$ws.Range("J136").Value = 2300
$ws.Range("K136").Value = 2691
$ws.Range("L136").Value = 6900
$ws.Range("M136").Value = -141
$ws.Range("N136").Value = -12000
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H104").Value = 3273.111
$ws.Range("J104").Value = 4065.4285
$ws.Range("L104").Value = 12196.2855
$ws.Range("N104").Value = -17438.2855
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 28577022
$ws.Range("I132").Value = 58831700
$ws.Range("J132").Value = 3158.6667
$ws.Range("K132").Value = 176495100
$ws.Range("L132").Value = 9476.000100000001
$ws.Range("M132").Value = -176492570
$ws.Range("N132").Value = -14536.0001
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4268.8887
$ws.Range("I40").Value = 3783
$ws.Range("K40").Value = 3783
$ws.Range("M40").Value = -3647
$ws.Range("H122").Value = 4719.7
$ws.Range("I122").Value = 3091
$ws.Range("K122").Value = 9273
$ws.Range("M122").Value = -6823
$ws.Range("H132").Value = 3320.16
$ws.Range("I132").Value = 2923.2307
$ws.Range("K132").Value = 8769.6921
$ws.Range("M132").Value = -6239.6921
$ws.Range("H135").Value = 52135
$ws.Range("J135").Value = 52135
$ws.Range("L135").Value = 52135
$ws.Range("N135").Value = -62275
$ws.Range("H136").Value = 1623.6786
$ws.Range("I136").Value = 1338.92
$ws.Range("J136").Value = 3996.6667
$ws.Range("K136").Value = 4016.76
$ws.Range("L136").Value = 11990.0001
$ws.Range("M136").Value = -1466.76
$ws.Range("N136").Value = -17090.0001
$ws.Range("H140").Value = 61461.777
$ws.Range("J140").Value = 61461.777
$ws.Range("L140").Value = 61461.777
$ws.Range("N140").Value = -71821.777
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 6208006
$ws.Range("I132").Value = 2370.2334
$ws.Range("K132").Value = 7110.7002
$ws.Range("M132").Value = -4580.7002
$ws.Range("H133").Value = 53428.75
$ws.Range("J133").Value = 53428.75
$ws.Range("L133").Value = 53428.75
$ws.Range("N133").Value = -63548.75
$ws.Range("H136").Value = 5167.6
$ws.Range("I136").Value = 5757.7144
$ws.Range("J136").Value = 4651.25
$ws.Range("K136").Value = 17273.1432
$ws.Range("L136").Value = 13953.75
$ws.Range("M136").Value = -14723.1432
$ws.Range("N136").Value = -19053.75
